# Add 2022-Q3 data
#
# 1) Insert a new worksheet named "2022-Q3" right after "总计" and before
#    "2022-Q2" (i.e. at position 2), and fill it with the quarterly fund
#    holdings detail.
# 2) Update the "总计" (summary) sheet: insert a new row for "2022-Q3" at
#    the top of the data (row 2) and push the existing quarters down by
#    one row, appending the previously-last "2020-Q4" row as a new row 9.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # "总计" summary sheet
$q2Sheet = $wb.Worksheets.Item(2)      # current "2022-Q2" sheet (about to become 3rd)

# ---------------------------------------------------------------------
# Step 1: create & populate the new "2022-Q3" detail sheet
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 2   # headers start at column B
    $cell = $newSheet.Cells.Item(1, $col)
    $ws1.Range("B1").Copy($cell)      # pick up the bold header style (s=2)
    $cell.Value = $headers[$i]
}

$data = @(
    @("014179", "中银证券远见价值混合A", "1.56", "93.65", "4.79", "0.0747", 5),
    @("005571", "中银证券新能源灵活配置混合A", "0.53", "90.32", "9.50", "0.0504", 1),
    @("003980", "中银证券瑞益灵活配置混合A", "0.66", "91.21", "6.38", "0.0421", 2),
    @("162216", "泰达宏利中证500指数增强（LOF）", "3.05", "93.17", "1.33", "0.0406", 6),
    @("005572", "中银证券新能源灵活配置混合C", "0.25", "90.32", "9.50", "0.0238", 1),
    @("003981", "中银证券瑞益灵活配置混合C", "0.19", "91.21", "6.38", "0.0121", 2),
    @("014180", "中银证券远见价值混合C", "0.16", "93.65", "4.79", "0.0077", 5),
    @("006783", "红土创新中证500指数增强A", "0.26", "92.93", "2.29", "0.0060", 6),
    @("159620", "华夏中证智选500成长创新策略ETF", "0.36", "91.92", "1.47", "0.0053", 6),
    @("501069", "华宝标普中国Ａ股质量价值指数（LOF）", "0.14", "93.78", "2.94", "0.0041", 7),
    @("660011", "农银中证500指数", "0.54", "94.17", "0.53", "0.0029", 9),
    @("006784", "红土创新中证500指数增强C", "0.12", "92.93", "2.29", "0.0027", 6),
    @("012926", "民生加银中证500指数增强A", "0.20", "86.58", "0.97", "0.0019", 10),
    @("012927", "民生加银中证500指数增强C", "0.14", "86.58", "0.97", "0.0014", 10)
)

# Force text storage (keeps leading zeros in fund codes, matches the
# "numeric-looking string" columns used throughout this workbook).
$newSheet.Range("B2:B15").NumberFormat = "@"
$newSheet.Range("D2:G15").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]

    $cellA = $newSheet.Cells.Item($row, 1)
    $ws1.Range("A2").Copy($cellA)     # pick up bold style (s=2) used on col A
    $cellA.Value = $i

    $newSheet.Cells.Item($row, 2).Value = $rowData[0]
    $newSheet.Cells.Item($row, 3).Value = $rowData[1]
    $newSheet.Cells.Item($row, 4).Value = $rowData[2]
    $newSheet.Cells.Item($row, 5).Value = $rowData[3]
    $newSheet.Cells.Item($row, 6).Value = $rowData[4]
    $newSheet.Cells.Item($row, 7).Value = $rowData[5]
    $newSheet.Cells.Item($row, 8).Value = $rowData[6]
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------
# Append a new row 9, copying the styled A8:D8 row so the new A9 cell
# keeps the bold "index column" style.
$ws1.Range("A8:D8").Copy($ws1.Range("A9:D9"))

# Shift the quarterly rows down by one (bottom-up so nothing is lost).
$ws1.Range("B9").Value = $ws1.Range("B8").Value()
$ws1.Range("C9").Value = $ws1.Range("C8").Value()
$ws1.Range("D9").Value = $ws1.Range("D8").Value()
$ws1.Range("A9").Value = 7

$ws1.Range("B8").Value = $ws1.Range("B7").Value()
$ws1.Range("C8").Value = $ws1.Range("C7").Value()
$ws1.Range("D8").Value = $ws1.Range("D7").Value()

$ws1.Range("B7").Value = $ws1.Range("B6").Value()
$ws1.Range("C7").Value = $ws1.Range("C6").Value()
$ws1.Range("D7").Value = $ws1.Range("D6").Value()

$ws1.Range("B6").Value = $ws1.Range("B5").Value()
$ws1.Range("C6").Value = $ws1.Range("C5").Value()
$ws1.Range("D6").Value = $ws1.Range("D5").Value()

$ws1.Range("B5").Value = $ws1.Range("B4").Value()
$ws1.Range("C5").Value = $ws1.Range("C4").Value()
$ws1.Range("D5").Value = $ws1.Range("D4").Value()

$ws1.Range("B4").Value = $ws1.Range("B3").Value()
$ws1.Range("C4").Value = $ws1.Range("C3").Value()
$ws1.Range("D4").Value = $ws1.Range("D3").Value()

$ws1.Range("B3").Value = $ws1.Range("B2").Value()
$ws1.Range("C3").Value = $ws1.Range("C2").Value()
$ws1.Range("D3").Value = $ws1.Range("D2").Value()

# New first data row: 2022-Q3
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 14
$ws1.Range("D2").Value = 0.28
